# Natmi following Dr Hou advice
# Update the LR-pair table (rows 2-4) with new values and append new rows 5-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array holds columns A..T (20 values) for sheet rows 2..7.
$rowsData = @(
  @("ECs", "Myl9", "Cd69", "ECs", 3, 1, 2.032726333333333, 6.098179, 0.006853963198393305, 0.006853963198393305, 3, 1, 36.51526533333334, 109.545796, 0.9971632597677756, 0.9971632597677755, 74.22554141172046, 668.0298727054841, 0.006834520285238238, 0.006834520285238238),
  @("ECs", "Myl9", "Cd69", "sCs", 3, 1, 2.032726333333333, 6.098179, 0.006853963198393305, 0.006853963198393305, 2, 0.6666666666666666, 0.103879, 0.311637, 0.002836740232224432, 0.002836740232224432, 0.2111575787803333, 1.900418209023, 0.00001944291315506794, 0.00001944291315506794),
  @("FAPs", "Myl9", "Cd69", "ECs", 3, 1, 4.690975, 14.072925, 0.01581706769246182, 0.01581706769246182, 3, 1, 36.51526533333334, 109.545796, 0.9971632597677756, 0.9971632597677755, 171.2921967970334, 1541.6297711733, 0.01577219878018279, 0.01577219878018279),
  @("FAPs", "Myl9", "Cd69", "sCs", 3, 1, 4.690975, 14.072925, 0.01581706769246182, 0.01581706769246182, 2, 0.6666666666666666, 0.103879, 0.311637, 0.002836740232224432, 0.002836740232224432, 0.487293792025, 4.385644128225, 0.0000448689122790237, 0.0000448689122790237),
  @("sCs", "Myl9", "Cd69", "ECs", 3, 1, 289.8530783333333, 869.559235, 0.9773289691091449, 0.9773289691091449, 3, 1, 36.51526533333334, 109.545796, 0.9971632597677756, 0.9971632597677755, 10584.06206302512, 95256.55856722606, 0.9745565407023546, 0.9745565407023545),
  @("sCs", "Myl9", "Cd69", "sCs", 3, 1, 289.8530783333333, 869.559235, 0.9773289691091449, 0.9773289691091449, 2, 0.6666666666666666, 0.103879, 0.311637, 0.002836740232224432, 0.002836740232224432, 30.10964792418833, 270.986831317695, 0.002772428406790341, 0.002772428406790341)
)

$r = 2
foreach ($rowVals in $rowsData) {
  $c = 1
  foreach ($val in $rowVals) {
    $ws.Cells.Item($r, $c).Value = $val
    $c = $c + 1
  }
  $r = $r + 1
}
